$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update FR-H03 description text (row 4, column C)
$ws.Range("C4").Value = 'The system shall allow the user to view the suggested destinations and remove unwanted stops or add new stops, after which the system shall automatically regenerate the optimized route.'

# 2. Make room for new rows by inserting blank rows at the right spots
# 2a. 4 new Functional rows (FR-L08..FR-L11) before the EIR-H01 row (row 9)
$ws.Rows("9:12").Insert()
# 2b. 4 new External Interface rows (EIR-L02..EIR-L05) after EIR-H01 (now row 13)
$ws.Rows("14:17").Insert()
# 2c. 2 new Performance rows (PR-L03..PR-L04) after PR-H02 (now row 19)
$ws.Rows("20:21").Insert()

# 3. Fill in the new Functional rows (9-12)
$ws.Range("A9").Value = 'FR-L08'
$ws.Range("B9").Value = 'Functional'
$ws.Range("C9").Value = 'The system shall automatically regenerate the optimized route when a user removes or adds a destination within the itinerary view.'
$ws.Range("A10").Value = 'FR-L09'
$ws.Range("B10").Value = 'Functional'
$ws.Range("C10").Value = 'The user shall be able to repeatedly modify the itinerary (e.g., remove or add POIs), and the system shall automatically regenerate the optimized route accordingly.'
$ws.Range("A11").Value = 'FR-L10'
$ws.Range("B11").Value = 'Functional'
$ws.Range("C11").Value = 'The system shall provide route-generation and recommendation support for at least twenty (20) cities or regions across Türkiye.'
$ws.Range("A12").Value = 'FR-L11'
$ws.Range("B12").Value = 'Functional'
$ws.Range("C12").Value = 'The system shall present a brief summary including total estimated distance and duration once a new route is generated.'

# 4. Fill in the new External Interface rows (14-17)
$ws.Range("A14").Value = 'EIR-L02'
$ws.Range("B14").Value = 'External Interface'
$ws.Range("C14").Value = 'The system shall allow the user to click on pins (markers) displayed on the map to open a detailed information page for the selected Point of Interest (POI).'
$ws.Range("A15").Value = 'EIR-L03'
$ws.Range("B15").Value = 'External Interface'
$ws.Range("C15").Value = 'The system shall communicate with OpenStreetMap layers through standardized API requests compatible with OSRM data format'
$ws.Range("A16").Value = 'EIR-L04'
$ws.Range("B16").Value = 'External Interface'
$ws.Range("C16").Value = 'The system shall maintain bidirectional navigation between the map interface and POI detail page, allowing the user to return to the map view with one action.'
$ws.Range("A17").Value = 'EIR-L05'
$ws.Range("B17").Value = 'External Interface'
$ws.Range("C17").Value = 'The system shall validate user inputs and display a clear error message without performing route calculation.'

# 5. Fill in the new Performance rows (20-21)
$ws.Range("A20").Value = 'PR-L03'
$ws.Range("B20").Value = 'Performance'
$ws.Range("C20").Value = 'The system shall start audio playback within 2 seconds after the user requests TTS for a POI description.'
$ws.Range("A21").Value = 'PR-L04'
$ws.Range("B21").Value = 'Performance'
$ws.Range("C21").Value = 'The system shall cache static POI and media data locally to minimize repeated API calls and improve response time.'

# 6. Remove the old trailing placeholder row (originally row 14, now shifted to row 24)
$ws.Rows("24:24").Delete()

# 7. Extend the table range to cover the new rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C21"))

# 8. Widen column C slightly to fit the new content
$ws.Columns("C").ColumnWidth = 142.0

# 9. Conditional formatting: drop the stray C14 anchor, keep only the header rule
$fc = $ws.Range("A1:C1").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("A1:C1"))

Write-Host "Edit complete"
